$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expansion List")

# Copy style (format) of the existing data row 13 down into new rows 14-16
$ws.Range("A13:F13").Copy() | Out-Null
$ws.Range("A14:F16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 14: A=10901-7(existing) D=2021-09(new) B=Display for 2021-09(new) C=SNOMEDCT(existing) E=...96(existing) F=FN(existing)
$ws.Range("A14").Value = "10901-7"
$ws.Range("D14").Value = "2021-09"
$ws.Range("B14").Value = "Display for 2021-09"
$ws.Range("C14").Value = "SNOMEDCT"
$ws.Range("E14").Value = "2.16.840.1.113883.6.96"
$ws.Range("F14").Value = "FN"

# Row 15: A=10901-8(new) B=Display for 10901-8(new) C=SNOMEDCT D=2021-09 E=...96 F=FN
$ws.Range("A15").Value = "10901-8"
$ws.Range("B15").Value = "Display for 10901-8"
$ws.Range("C15").Value = "SNOMEDCT"
$ws.Range("D15").Value = "2021-09"
$ws.Range("E15").Value = "2.16.840.1.113883.6.96"
$ws.Range("F15").Value = "FN"

# Row 16: A=10901-8(existing) C=LOINC(new) B=Display for 10901-8 LOINC(new) D=2021-09(existing) E=2.16.840.1.113883.6.1(new) F=FN(existing)
$ws.Range("A16").Value = "10901-8"
$ws.Range("C16").Value = "LOINC"
$ws.Range("B16").Value = "Display for 10901-8 LOINC"
$ws.Range("D16").Value = "2021-09"
$ws.Range("E16").Value = "2.16.840.1.113883.6.1"
$ws.Range("F16").Value = "FN"

# Row 16 has a custom height set (14.25)
$ws.Rows.Item(16).RowHeight = 14.25

# Update selection/active cell to match target
$ws.Range("A16:XFD16").Select() | Out-Null
